# Convert the CONCATENATE formulas in column L (rows 2-31) of Hoja2 into
# their static computed text values (formula -> value), equivalent to a
# "Copy" + "Paste Special: Values" over L2:L31, and move the selection there.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 'INSERT INTO ''persona'' VALUES ("Usmar Isacc","Miguel","Lopez",16/10/1999,"M","LopezMiguelUs@gmail.com","9717995120","direccion prueba","Ixtepec","/directory/img-person/UsmarIsaccMiguel.jpg"'
$ws.Range("L3").Value = 'INSERT INTO ''persona'' VALUES ("Giovanni Alexander","Morales","Lopez",11/04/1999,"M","LopezMoralesGi@gmail.com","9714324512","direccion prueba","Ixtaltepec","/directory/img-person/GiovanniAlexanderMorales.jpg"'
$ws.Range("L4").Value = 'INSERT INTO ''persona'' VALUES ("Roel","Morales","Orozco",06/05/1998,"M","OrozcoMoralesRo@gmail.com","9714860957","direccion prueba","Juchitán","/directory/img-person/RoelMorales.jpg"'
$ws.Range("L5").Value = 'INSERT INTO ''persona'' VALUES ("Karina","Nuñez","Silva",07/05/1999,"F","SilvaNuñezKa@gmail.com","9712623032","direccion prueba","Ixtepec","/directory/img-person/KarinaNuñez.jpg"'
$ws.Range("L6").Value = 'INSERT INTO ''persona'' VALUES ("Kenia","Orozco","Cruz",31/10/1999,"F","CruzOrozcoKe@gmail.com","9712793311","direccion prueba","Salina Cruz","/directory/img-person/KeniaOrozco.jpg"'
$ws.Range("L7").Value = 'INSERT INTO ''persona'' VALUES ("Dennis Sabas","Orozco","Bautista",23/01/1999,"F","BautistaOrozcoDe@gmail.com","9713721570","direccion prueba","Juchitán","/directory/img-person/DennisSabasOrozco.jpg"'
$ws.Range("L8").Value = 'INSERT INTO ''persona'' VALUES ("Carlos Roberto","Ortiz","Matus",17/05/1999,"M","MatusOrtizCa@gmail.com","9718332712","direccion prueba","Salina Cruz","/directory/img-person/CarlosRobertoOrtiz.jpg"'
$ws.Range("L9").Value = 'INSERT INTO ''persona'' VALUES ("Isis Yamile","Altamirano","Solano",03/05/1998,"F","SolanoAltamiranoIs@gmail.com","9718259758","direccion prueba","Ixtepec","/directory/img-person/IsisYamileAltamirano.jpg"'
$ws.Range("L10").Value = 'INSERT INTO ''persona'' VALUES ("Fátima","Antinio","Sanchez",08/01/1999,"F","SanchezAntinioFá@gmail.com","9716330120","direccion prueba","Juchitán","/directory/img-person/FátimaAntinio.jpg"'
$ws.Range("L11").Value = 'INSERT INTO ''persona'' VALUES ("Itzury Alejandra","Aquino","Cruz",27/12/1999,"F","CruzAquinoIt@gmail.com","9716802586","direccion prueba","Chicapa","/directory/img-person/ItzuryAlejandraAquino.jpg"'
$ws.Range("L12").Value = 'INSERT INTO ''persona'' VALUES ("Vicente","Aquino","Regalado",11/10/1999,"M","RegaladoAquinoVi@gmail.com","9712106115","direccion prueba","Juchitán","/directory/img-person/VicenteAquino.jpg"'
$ws.Range("L13").Value = 'INSERT INTO ''persona'' VALUES ("Baudel","Aranjo","Benitez",03/12/1999,"M","BenitezAranjoBa@gmail.com","9717592062","direccion prueba","Espinal","/directory/img-person/BaudelAranjo.jpg"'
$ws.Range("L14").Value = 'INSERT INTO ''persona'' VALUES ("Yosmar Manuel","Avedaño","Morales",14/11/1999,"M","MoralesAvedañoYo@gmail.com","9716019661","direccion prueba","Ixtaltepec","/directory/img-person/YosmarManuelAvedaño.jpg"'
$ws.Range("L15").Value = 'INSERT INTO ''persona'' VALUES ("Brian","Benefield","Morales",26/07/1999,"M","MoralesBenefieldBr@gmail.com","9714772466","direccion prueba","Unión Hidalgo","/directory/img-person/BrianBenefield.jpg"'
$ws.Range("L16").Value = 'INSERT INTO ''persona'' VALUES ("Kevin Gabriel","Zarate","Velasquez",08/07/1998,"M","VelasquezZarateKe@gmail.com","9712498199","direccion prueba","Chicapa","/directory/img-person/KevinGabrielZarate.jpg"'
$ws.Range("L17").Value = 'INSERT INTO ''persona'' VALUES ("Jesus Antonio","Zarate","Villalobos",23/04/1999,"M","VillalobosZarateJe@gmail.com","9716213801","direccion prueba","Juchitán","/directory/img-person/JesusAntonioZarate.jpg"'
$ws.Range("L18").Value = 'INSERT INTO ''persona'' VALUES ("Clarissa","Zavala","Jiménez",01/09/1999,"F","JiménezZavalaCl@gmail.com","9715339317","direccion prueba","Chicapa","/directory/img-person/ClarissaZavala.jpg"'
$ws.Range("L19").Value = 'INSERT INTO ''persona'' VALUES ("Alba Beatriz","Aguilar","Ulises",03/06/1999,"F","UlisesAguilarAl@gmail.com","9714414298","direccion prueba","Ixtepec","/directory/img-person/AlbaBeatrizAguilar.jpg"'
$ws.Range("L20").Value = 'INSERT INTO ''persona'' VALUES ("Luis","Fuentes","de Jesus",21/12/1999,"M","de JesusFuentesLu@gmail.com","9713345356","direccion prueba","Espinal","/directory/img-person/LuisFuentes.jpg"'
$ws.Range("L21").Value = 'INSERT INTO ''persona'' VALUES ("Maytor","Revuelta","Rosado",27/11/1999,"M","RosadoRevueltaMa@gmail.com","9713195036","direccion prueba","Juchitán","/directory/img-person/MaytorRevuelta.jpg"'
$ws.Range("L22").Value = 'INSERT INTO ''persona'' VALUES ("Jorge","Robles","Luis",27/04/1999,"M","LuisRoblesJo@gmail.com","9714341427","direccion prueba","Tehuantepec","/directory/img-person/JorgeRobles.jpg"'
$ws.Range("L23").Value = 'INSERT INTO ''persona'' VALUES ("Sergio","Rodas","Escobar",27/11/1998,"M","EscobarRodasSe@gmail.com","9713849923","direccion prueba","Unión Hidalgo","/directory/img-person/SergioRodas.jpg"'
$ws.Range("L24").Value = 'INSERT INTO ''persona'' VALUES ("Brisa Donaji","Ruiz","Sanchez",08/01/1999,"F","SanchezRuizBr@gmail.com","9716521066","direccion prueba","Ixtaltepec","/directory/img-person/BrisaDonajiRuiz.jpg"'
$ws.Range("L25").Value = 'INSERT INTO ''persona'' VALUES ("Jair Michael","Ruiz","Vicente",12/10/1999,"M","VicenteRuizJa@gmail.com","9715469345","direccion prueba","Tehuantepec","/directory/img-person/JairMichaelRuiz.jpg"'
$ws.Range("L26").Value = 'INSERT INTO ''persona'' VALUES ("José Guillermo","Gonzales","Lopez",22/01/1999,"M","LopezGonzalesJo@gmail.com","9711561046","direccion prueba","Juchitán","/directory/img-person/JoséGuillermoGonzales.jpg"'
$ws.Range("L27").Value = 'INSERT INTO ''persona'' VALUES ("Francisco","Santiago","de la Cruz",29/10/1998,"M","CruzSantiagoFr@gmail.com","9711744464","direccion prueba","Juchitán","/directory/img-person/FranciscoSantiago.jpg"'
$ws.Range("L28").Value = 'INSERT INTO ''persona'' VALUES ("Luis Alberto","Robles","Parada",03/04/1998,"M","ParadaRoblesLu@gmail.com","9719614394","direccion prueba","Tehuantepec","/directory/img-person/Luis AlbertoRobles.jpg"'
$ws.Range("L29").Value = 'INSERT INTO ''persona'' VALUES ("Juan Carlos","Fernández","Piñon",03/04/1998,"M","PiñonFernándezJu@gmail.com","9715787642","direccion prueba","Ixtepec","/directory/img-person/Juan CarlosFernández.jpg"'
$ws.Range("L30").Value = 'INSERT INTO ''persona'' VALUES ("Emanuel","Enríquez","Couder",30/11/1998,"M","CouderEnríquezEm@gmail.com","9719665488","direccion prueba","Guevea","/directory/img-person/EmanuelEnríquez.jpg"'
$ws.Range("L31").Value = 'INSERT INTO ''persona'' VALUES ("Humberto","Toledo","Fuentes",23/06/1998,"M","FuentesToledoHu@gmail.com","9718507158","direccion prueba","Ixtepec","/directory/img-person/HumbertoToledo.jpg"'

$ws.Range("L2:L31").Select()
$wb.Save()
